$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(14).Insert()
$ws.Rows.Item(20).Insert()
$ws.Rows.Item(24).Insert()
$ws.Rows.Item(26).Insert()
$ws.Rows.Item(27).Insert()
$ws.Rows.Item(28).Insert()

for ($r = 2; $r -le 29; $r++) {
    $v = $ws.Range("A$r").Value2
    Write-Host ("$r A=[" + $v + "]")
}
